$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "statut" column (A) used four emoji values as status flags.
# They are replaced with plain-text/emoji markers that render more
# reliably across Excel builds/fonts:
#   old blue book   (U+1F4D8) -> warning sign "⚠️"
#   old red book    (U+1F4D5) -> "-3"
#   old orange book (U+1F4D9) -> "+3"
#   old green book  (U+1F4D7) -> check mark "✅"
$oldBlue   = "📘"
$oldRed    = "📕"
$oldOrange = "📙"
$oldGreen  = "📗"

$newWarning = "⚠️"
$newCheck   = "✅"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $used.Rows.Count + $used.Row - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -eq $null) { continue }
    $sval = [string]$val

    if ($sval -eq $oldBlue) {
        # Plain text value, not numeric-looking, safe to assign directly.
        $cell.Value = $newWarning
    }
    elseif ($sval -eq $oldGreen) {
        $cell.Value = $newCheck
    }
    elseif ($sval -eq $oldRed) {
        # "-3" looks like a number to Excel's auto-detection, which would
        # store it as a numeric cell instead of text. Build it through a
        # text-literal formula, then convert the formula to its resulting
        # value via copy / paste-special so the cell keeps its original
        # (default) style while still ending up as a shared text string.
        $cell.Formula = '="-3"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
    elseif ($sval -eq $oldOrange) {
        $cell.Formula = '="+3"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = 0
